$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("testdata")

$ws.Range("A9").Value = "ajay"
$ws.Range("A9").Select()
